# Mediciones.xlsx - TP5 - "mediciones de les cuadripepes"
#
# Updates the two "Cuadripolo" measurement tables on Hoja1:
#   - Renames/updates the device id strings (A1 / A8)
#   - Fills in the measured V/I values for both tables (rows 3-6 and 10-13)
#   - Applies a 2-decimal number format to the three hand-measured Z values
#   - Updates a couple of downstream ratio formulas in row 11
#   - Restores the previous selection used on Hoja1

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")

# --- Device identification labels -----------------------------------------
$ws1.Range("A1").Value2 = 'Cuadripolo "A" - Nro. De serie: 9608 (15 V - 50 mA) (T)'
$ws1.Range("A8").Value2 = 'Cuadripolo "B" - Nro. De serie: 9603 (15 V - 50 mA) (Pi)'

# --- Cuadripolo "A" table (rows 3-6) ---------------------------------------
$ws1.Range("C3").Value2 = 3.05
$ws1.Range("D3").Formula = "=0.5*23"
$ws1.Range("E3").Formula = "=32*0.5"

$ws1.Range("B4").Value2 = 3.05
$ws1.Range("D4").Formula = "=33.5*0.5"
$ws1.Range("E4").Formula = "=0.5*23"

$ws1.Range("B5").Formula = "=2.11"
$ws1.Range("C5").Value2 = 3.05
$ws1.Range("E5").Formula = "=0.1*68"

$ws1.Range("B6").Value2 = 3.05
$ws1.Range("C6").NumberFormat = "0.00"
$ws1.Range("C6").Formula = "=2.185"
$ws1.Range("D6").Formula = "=82*0.1"

# --- Cuadripolo "B" table (rows 10-13) -------------------------------------
$ws1.Range("C10").Value2 = 3.05
$ws1.Range("D10").Formula = "=83*0.5"
$ws1.Range("E10").Formula = "=60*0.5"

$ws1.Range("B11").Value2 = 3.05
$ws1.Range("D11").Formula = "=72*0.5"
$ws1.Range("E11").Formula = "=60.5*0.5"

$ws1.Range("B12").NumberFormat = "0.00"
$ws1.Range("B12").Formula = "=2.557"
$ws1.Range("C12").Formula = "=3.05"
$ws1.Range("E12").Formula = "=32.5*0.5"

$ws1.Range("B13").Value2 = 3.05
$ws1.Range("C13").NumberFormat = "0.00"
$ws1.Range("C13").Formula = "=2.203"
$ws1.Range("D13").Formula = "=28*0.5"

# Row 11 downstream ratio formulas were edited by hand alongside the new
# data and now reference the adjusted cells (two of them now point at a
# removed column, matching the author's edit).
$ws1.Range("H11").Formula = "=D13/#REF!*1000"
$ws1.Range("J11").Formula = "=C13/#REF!*1000"
$ws1.Range("R11").Formula = "=D13/C13"
$ws1.Range("T11").Formula = "=#REF!/C13/1000"

# --- Restore the selection left on Hoja1 ------------------------------------
$ws1.Range("G9").Select()
